$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 666.6667
$ws.Range("J2").Value = 950
$ws.Range("L2").Value = 950
$ws.Range("N2").Value = -1176
$ws.Range("H8").Value = 29.285715
$ws.Range("I8").Value = 29.285715
$ws.Range("K8").Value = 87.857145
$ws.Range("M8").Value = 51.142855
$ws.Range("H62").Value = 13330.889
$ws.Range("I62").Value = 10829.667
$ws.Range("K62").Value = 10829.667
$ws.Range("M62").Value = -10205.667
$ws.Range("H65").Value = 13330.889
$ws.Range("I65").Value = 10829.667
$ws.Range("K65").Value = 54148.335
$ws.Range("M65").Value = -51028.335
$ws.Range("H98").Value = 973.30554
$ws.Range("I98").Value = 1019.6667
$ws.Range("K98").Value = 1019.6667
$ws.Range("M98").Value = 478.3333
$ws.Range("H112").Value = 2650.853
$ws.Range("I112").Value = 2159.6667
$ws.Range("K112").Value = 6479.000100000001
$ws.Range("M112").Value = -5371.000100000001
$ws.Range("H122").Value = 973.30554
$ws.Range("I122").Value = 1019.6667
$ws.Range("K122").Value = 3059.0001
$ws.Range("M122").Value = -609.0001000000002
$ws.Range("H131").Value = 13978.357
$ws.Range("I131").Value = 1376.1
$ws.Range("J131").Value = 45484
$ws.Range("K131").Value = 4128.299999999999
$ws.Range("L131").Value = 136452
$ws.Range("M131").Value = 911.7000000000007
$ws.Range("N131").Value = -146532
$ws.Range("H135").Value = 4608.7334
$ws.Range("I135").Value = 813.1
$ws.Range("K135").Value = 7317.900000000001
$ws.Range("M135").Value = -4782.900000000001
$ws.Range("H137").Value = 27781264
$ws.Range("I137").Value = 58824556
$ws.Range("K137").Value = 176473668
$ws.Range("M137").Value = -176471118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4004391.8
$ws.Range("I61").Value = 4574.8335
$ws.Range("K61").Value = 4574.8335
$ws.Range("M61").Value = -4362.8335
$ws.Range("H81").Value = 73639.664
$ws.Range("J81").Value = 73639.664
$ws.Range("L81").Value = 73639.664
$ws.Range("N81").Value = -75635.664
$ws.Range("H84").Value = 73639.664
$ws.Range("J84").Value = 73639.664
$ws.Range("L84").Value = 220918.992
$ws.Range("N84").Value = -230902.992
$ws.Range("H94").Value = 44330
$ws.Range("J94").Value = 44330
$ws.Range("L94").Value = 44330
$ws.Range("N94").Value = -46132
$ws.Range("H97").Value = 1680.5834
$ws.Range("I97").Value = 1349.238
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 1349.238
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = -853.2380000000001
$ws.Range("N97").Value = -4992
$ws.Range("H132").Value = 473656.66
$ws.Range("I132").Value = 533774.2
$ws.Range("J132").Value = 2735.8333
$ws.Range("K132").Value = 1601322.6
$ws.Range("L132").Value = 8207.499899999999
$ws.Range("M132").Value = -1598792.6
$ws.Range("N132").Value = -13267.4999
$ws.Range("H136").Value = 4004391.8
$ws.Range("I136").Value = 4574.8335
$ws.Range("K136").Value = 13724.5005
$ws.Range("M136").Value = -11174.5005
$ws.Range("H140").Value = 70349.5
$ws.Range("J140").Value = 70349.5
$ws.Range("L140").Value = 70349.5
$ws.Range("N140").Value = -80709.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21216.846
$ws.Range("I82").Value = 10574.333
$ws.Range("J82").Value = 45162.5
$ws.Range("K82").Value = 10574.333
$ws.Range("L82").Value = 45162.5
$ws.Range("M82").Value = -10191.333
$ws.Range("N82").Value = -45928.5
$ws.Range("H85").Value = 21216.846
$ws.Range("I85").Value = 10574.333
$ws.Range("J85").Value = 45162.5
$ws.Range("K85").Value = 10574.333
$ws.Range("L85").Value = 45162.5
$ws.Range("M85").Value = -9248.333000000001
$ws.Range("N85").Value = -47814.5
$ws.Range("H106").Value = 70000
$ws.Range("J106").Value = 70000
$ws.Range("L106").Value = 70000
$ws.Range("N106").Value = -72524
$ws.Range("H107").Value = 1367.8889
$ws.Range("I107").Value = 976.375
$ws.Range("J107").Value = 4500
$ws.Range("K107").Value = 976.375
$ws.Range("L107").Value = 4500
$ws.Range("M107").Value = 943.625
$ws.Range("N107").Value = -8340
$ws.Range("H134").Value = 6178460.5
$ws.Range("I134").Value = 6513.316
$ws.Range("K134").Value = 19539.948
$ws.Range("M134").Value = -17004.948

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 71482.57000000001
$ws.Range("J7").Value = 78
$ws.Range("L7").Value = 78
$ws.Range("N7").Value = -304
$ws.Range("H31").Value = 1426946.6
$ws.Range("I31").Value = 1738608.9
$ws.Range("J31").Value = 2204.7144
$ws.Range("K31").Value = 1738608.9
$ws.Range("L31").Value = 2204.7144
$ws.Range("M31").Value = -1738313.9
$ws.Range("N31").Value = -2794.7144
$ws.Range("H34").Value = 1426946.6
$ws.Range("I34").Value = 1738608.9
$ws.Range("J34").Value = 2204.7144
$ws.Range("K34").Value = 1738608.9
$ws.Range("L34").Value = 2204.7144
$ws.Range("M34").Value = -1738406.9
$ws.Range("N34").Value = -2608.7144
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H107").Value = 254.08
$ws.Range("I107").Value = 254.08
$ws.Range("K107").Value = 254.08
$ws.Range("M107").Value = 1665.92
$ws.Range("H132").Value = 2838.9512
$ws.Range("I132").Value = 2735.5264
$ws.Range("K132").Value = 8206.5792
$ws.Range("M132").Value = -5676.5792
$ws.Range("H135").Value = 71014.8
$ws.Range("J135").Value = 71014.8
$ws.Range("L135").Value = 71014.8
$ws.Range("N135").Value = -81154.8
$ws.Range("H141").Value = 168452.94
$ws.Range("J141").Value = 199733.75
$ws.Range("L141").Value = 199733.75
$ws.Range("N141").Value = -210093.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2658
$ws.Range("I70").Value = 3249.5
$ws.Range("J70").Value = 1475
$ws.Range("K70").Value = 9748.5
$ws.Range("L70").Value = 4425
$ws.Range("M70").Value = -9433.5
$ws.Range("N70").Value = -5055
$ws.Range("H73").Value = 2658
$ws.Range("I73").Value = 3249.5
$ws.Range("J73").Value = 1475
$ws.Range("K73").Value = 9748.5
$ws.Range("L73").Value = 4425
$ws.Range("M73").Value = -8656.5
$ws.Range("N73").Value = -6609

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 40082.64
$ws.Range("I2").Value = 55611.945
$ws.Range("K2").Value = 55611.945
$ws.Range("M2").Value = -55498.945
$ws.Range("H26").Value = 15000
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H50").Value = 15000
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H70").Value = 49472
$ws.Range("I70").Value = 13114.538
$ws.Range("K70").Value = 13114.538
$ws.Range("M70").Value = -12844.538
$ws.Range("H73").Value = 49472
$ws.Range("I73").Value = 13114.538
$ws.Range("K73").Value = 13114.538
$ws.Range("M73").Value = -12178.538
$ws.Range("H97").Value = 1104.68
$ws.Range("I97").Value = 1115.8667
$ws.Range("K97").Value = 1115.8667
$ws.Range("M97").Value = -619.8667
$ws.Range("H102").Value = 1721.28
$ws.Range("I102").Value = 1668
$ws.Range("K102").Value = 1668
$ws.Range("M102").Value = -46
$ws.Range("H113").Value = 3206.5833
$ws.Range("I113").Value = 2813.4285
$ws.Range("K113").Value = 2813.4285
$ws.Range("M113").Value = -643.4285
$ws.Range("H122").Value = 92966
$ws.Range("I122").Value = 134261.62
$ws.Range("J122").Value = 10374.75
$ws.Range("K122").Value = 402784.86
$ws.Range("L122").Value = 31124.25
$ws.Range("M122").Value = -400334.86
$ws.Range("N122").Value = -36024.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1265
$ws.Range("I16").Value = 328.16666
$ws.Range("K16").Value = 328.16666
$ws.Range("M16").Value = -158.16666
$ws.Range("H56").Value = 5743.3335
$ws.Range("I56").Value = 5743.3335
$ws.Range("K56").Value = 5743.3335
$ws.Range("M56").Value = -5052.3335
$ws.Range("H74").Value = 49671.5
$ws.Range("I74").Value = 49671.5
$ws.Range("K74").Value = 49671.5
$ws.Range("M74").Value = -48673.5
$ws.Range("H77").Value = 49671.5
$ws.Range("I77").Value = 49671.5
$ws.Range("K77").Value = 149014.5
$ws.Range("M77").Value = -144022.5
$ws.Range("H82").Value = 2140.889
$ws.Range("I82").Value = 3267
$ws.Range("J82").Value = 1240
$ws.Range("K82").Value = 3267
$ws.Range("L82").Value = 1240
$ws.Range("M82").Value = -2906
$ws.Range("N82").Value = -1962
$ws.Range("H85").Value = 2140.889
$ws.Range("I85").Value = 3267
$ws.Range("J85").Value = 1240
$ws.Range("K85").Value = 3267
$ws.Range("L85").Value = 1240
$ws.Range("M85").Value = -2019
$ws.Range("N85").Value = -3736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 13967.333
$ws.Range("I61").Value = 13967.333
$ws.Range("K61").Value = 13967.333
$ws.Range("M61").Value = -13675.333
$ws.Range("H126").Value = 2663.55
$ws.Range("I126").Value = 2593.2632
$ws.Range("K126").Value = 7779.7896
$ws.Range("M126").Value = -5309.7896
$ws.Range("H141").Value = 89995
$ws.Range("J141").Value = 89995
$ws.Range("L141").Value = 89995
$ws.Range("N141").Value = -100355
